$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:L (and the remaining default columns) are marked as
# outline-collapsed in the saved file.
$ws.Columns.Item(1).ShowDetail = $false
$ws.Columns.Item(2).ShowDetail = $false
$ws.Columns.Item(3).ShowDetail = $false
$ws.Columns.Item(4).ShowDetail = $false
$ws.Columns.Item(5).ShowDetail = $false
$ws.Columns.Item(6).ShowDetail = $false
$ws.Columns.Item(7).ShowDetail = $false
$ws.Columns.Item(8).ShowDetail = $false
$ws.Columns.Item(9).ShowDetail = $false
$ws.Columns.Item(10).ShowDetail = $false
$ws.Columns.Item(11).ShowDetail = $false

# D10 changes from 21 to 100 (numeric).
$ws.Range("D10").Value = 100
